# Skill.xlsx edit — "Added size for the scene / Fixed Datalist ToString bug /
# Added protocol for serializing the player position".
#
# The concrete, data-level change in the sheet is a balance tweak: every
# RequireDistance (Q) / DamageDistance (R) value for the skill rows (11-61)
# is bumped from 1.5/2 to a uniform 2.5. Along the way a handful of cells in
# column Q (rows 53-58) pick up the row-banding fill that neighbouring rows
# (50-52) already use, which we reproduce by copying the format across.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Activate()

# --- Data edit: RequireDistance / DamageDistance columns, rows 11-61 -------
# Q and R are adjacent, so one rectangular assignment covers every row.
$ws.Range("Q11:R61").Value = 2.5

# --- Formatting follow-on: Q53:Q58 adopt the same banded fill as Q50:Q52 ---
$ws.Range("Q50").Copy()
$ws.Range("Q53:Q58").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = 0

# --- Cosmetic view state (best effort) --------------------------------------
# The author also scrolled/resized the window while editing down near row 56
# and left the cursor on Q59; reproduce the parts of that view state the
# object model exposes.
$win = $excel.ActiveWindow
$win.ScrollColumn = 13
$win.ScrollRow = 56
$win.Height = 690

$ws.Range("Q59").Select()
